$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.535095978099406
$ws.Range("C2").Value = 0.3053878575661422
$ws.Range("D2").Value = 0.03678230788326431
$ws.Range("F2").Value = 0.7017207948563495
$ws.Range("G2").Value = 0.5430413825110492
$ws.Range("H2").Value = 0.6874466824658825
$ws.Range("K2").Value = 0.2870740952819233
$ws.Range("L2").Value = 0.3022841447732532
$ws.Range("O2").Value = 2.43268904452944
$ws.Range("B3").Value = 0.4877282317053755
$ws.Range("C3").Value = 0.308387659888929
$ws.Range("D3").Value = 0.03400285097652045
$ws.Range("F3").Value = 0.7035530904853431
$ws.Range("G3").Value = 0.5464289991049966
$ws.Range("H3").Value = 0.6927265321242544
$ws.Range("K3").Value = 0.2504400313880524
$ws.Range("L3").Value = 0.2910177986555738
$ws.Range("O3").Value = 2.450674007495664
$ws.Range("B4").Value = 0.4587191024306492
$ws.Range("C4").Value = 0.3103315796407902
$ws.Range("D4").Value = 0.03228202559675708
$ws.Range("F4").Value = 0.7051133395097082
$ws.Range("G4").Value = 0.5488822302595793
$ws.Range("H4").Value = 0.696265468618094
$ws.Range("K4").Value = 0.2278631197621905
$ws.Range("L4").Value = 0.2842715012239552
$ws.Range("O4").Value = 2.463121297605767
$ws.Range("B5").Value = 0.4469172379610029
$ws.Range("C5").Value = 0.311149426335815
$ws.Range("D5").Value = 0.03157723513504607
$ws.Range("F5").Value = 0.705858610222073
$ws.Range("G5").Value = 0.5499757283537576
$ws.Range("H5").Value = 0.6977823689010592
$ws.Range("K5").Value = 0.2186424485574463
$ws.Range("L5").Value = 0.281565476567323
$ws.Range("O5").Value = 2.468546794850823
$ws.Range("B6").Value = 0.4449587532498072
$ws.Range("C6").Value = 0.3112867811987599
$ws.Range("D6").Value = 0.03145999250524767
$ws.Range("F6").Value = 0.7059889733152147
$ws.Range("G6").Value = 0.5501629655007605
$ws.Range("H6").Value = 0.6980387655197475
$ws.Range("K6").Value = 0.2171101485676701
$ws.Range("L6").Value = 0.2811187513911619
$ws.Range("O6").Value = 2.469469021418661
$ws.Range("B7").Value = 0.4585598580156045
$ws.Range("C7").Value = 0.3103425053783671
$ws.Range("D7").Value = 0.03227253483058234
$ws.Range("F7").Value = 0.7051229472811329
$ws.Range("G7").Value = 0.5488965979032798
$ws.Range("H7").Value = 0.6962856233203212
$ws.Range("K7").Value = 0.2277388483695404
$ws.Range("L7").Value = 0.2842348320307622
$ws.Range("O7").Value = 2.463193038015703
$ws.Range("B8").Value = 0.5187484871656523
$ws.Range("C8").Value = 0.306401030085711
$ws.Range("D8").Value = 0.03582692578979874
$ws.Range("F8").Value = 0.7022622520610753
$ws.Range("G8").Value = 0.5441319322889129
$ws.Range("H8").Value = 0.6892055440355662
$ws.Range("K8").Value = 0.274460330549033
$ws.Range("L8").Value = 0.2983639996477478
$ws.Range("O8").Value = 2.43859874370915
$ws.Range("B9").Value = 0.6373444302820417
$ws.Range("C9").Value = 0.2994800407940925
$ws.Range("D9").Value = 0.04268284110201392
$ws.Range("F9").Value = 0.7001056071834029
$ws.Range("G9").Value = 0.5377530213948702
$ws.Range("H9").Value = 0.6776766607162514
$ws.Range("K9").Value = 0.3653974096094998
$ws.Range("L9").Value = 0.3274287793100257
$ws.Range("O9").Value = 2.401515786436093
$ws.Range("B10").Value = 0.7247935658799065
$ws.Range("C10").Value = 0.2948859683933609
$ws.Range("D10").Value = 0.04764887961316333
$ws.Range("F10").Value = 0.7006273785194708
$ws.Range("G10").Value = 0.5348787994104498
$ws.Range("H10").Value = 0.6706393831478863
$ws.Range("K10").Value = 0.4317696322953566
$ws.Range("L10").Value = 0.3496108035072751
$ws.Range("O10").Value = 2.381071429405011
$ws.Range("B11").Value = 0.7646397190200673
$ws.Range("C11").Value = 0.2929021616930481
$ws.Range("D11").Value = 0.04989237892871756
$ws.Range("F11").Value = 0.7013224054556559
$ws.Range("G11").Value = 0.5339659005428956
$ws.Range("H11").Value = 0.6677485368498992
$ws.Range("K11").Value = 0.4618643345237672
$ws.Range("L11").Value = 0.3598820693295437
$ws.Range("O11").Value = 2.3732486941324
$ws.Range("B12").Value = 0.7797371239623203
$ws.Range("C12").Value = 0.2921661696472828
$ws.Range("D12").Value = 0.05073966128000507
$ws.Range("F12").Value = 0.7016514087237766
$ws.Range("G12").Value = 0.5336770329001581
$ws.Range("H12").Value = 0.6666984465414032
$ws.Range("K12").Value = 0.4732457643102919
$ws.Range("L12").Value = 0.3637974545456188
$ws.Range("O12").Value = 2.370498976214122
$ws.Range("B13").Value = 0.7764852614394613
$ws.Range("C13").Value = 0.2923240015271578
$ws.Range("D13").Value = 0.05055728603797149
$ws.Range("F13").Value = 0.7015776250195245
$ws.Range("G13").Value = 0.5337367169709495
$ws.Range("H13").Value = 0.6669226187151622
$ws.Range("K13").Value = 0.4707952372711475
$ws.Range("L13").Value = 0.3629530566130796
$ws.Range("O13").Value = 2.371081720651688
$ws.Range("B14").Value = 0.765881625376494
$ws.Range("C14").Value = 0.2928413060088157
$ws.Range("D14").Value = 0.04996213133655658
$ws.Range("F14").Value = 0.7013481537983353
$ws.Range("G14").Value = 0.5339409958436079
$ws.Range("H14").Value = 0.6676612514078215
$ws.Range("K14").Value = 0.4628009914746372
$ws.Range("L14").Value = 0.3602036722400328
$ws.Range("O14").Value = 2.373018212049374
$ws.Range("B15").Value = 0.7593876785601594
$ws.Range("C15").Value = 0.2931601531106729
$ws.Range("D15").Value = 0.04959728344068992
$ws.Range("F15").Value = 0.7012161666227712
$ws.Range("G15").Value = 0.5340735256376234
$ws.Range("H15").Value = 0.6681194939263122
$ws.Range("K15").Value = 0.4579023413178902
$ws.Range("L15").Value = 0.3585229627428106
$ws.Range("O15").Value = 2.374232056902343
$ws.Range("B16").Value = 0.7221907650909429
$ws.Range("C16").Value = 0.2950177474066713
$ws.Range("D16").Value = 0.04750194460493162
$ws.Range("F16").Value = 0.7005911663359754
$ws.Range("G16").Value = 0.5349464070879577
$ws.Range("H16").Value = 0.670834551044905
$ws.Range("K16").Value = 0.4298008390275641
$ws.Range("L16").Value = 0.3489431781907655
$ws.Range("O16").Value = 2.381612404931815
$ws.Range("B17").Value = 0.6993877466242111
$ws.Range("C17").Value = 0.2961844723137013
$ws.Range("D17").Value = 0.04621250249965669
$ws.Range("F17").Value = 0.7003249691412563
$ws.Range("G17").Value = 0.535583018341022
$ws.Range("H17").Value = 0.6725796410905218
$ws.Range("K17").Value = 0.4125358321952035
$ws.Range("L17").Value = 0.3431124753045935
$ws.Range("O17").Value = 2.386518519088952
$ws.Range("B18").Value = 0.6862782308110127
$ws.Range("C18").Value = 0.2968655255770312
$ws.Range("D18").Value = 0.04546938589487581
$ws.Range("F18").Value = 0.7002149378028193
$ws.Range("G18").Value = 0.5359863168528989
$ws.Range("H18").Value = 0.6736125932975199
$ws.Range("K18").Value = 0.4025962372436709
$ws.Range("L18").Value = 0.3397758086576488
$ws.Range("O18").Value = 2.389479449015639
$ws.Range("B19").Value = 0.6818406602889411
$ws.Range("C19").Value = 0.2970978338420593
$ws.Range("D19").Value = 0.04521752929272793
$ws.Range("F19").Value = 0.7001850818171391
$ws.Range("G19").Value = 0.5361292424477
$ws.Range("H19").Value = 0.6739673531296759
$ws.Range("K19").Value = 0.3992292971576035
$ws.Range("L19").Value = 0.3386489928549707
$ws.Range("O19").Value = 2.390505849539068
$ws.Range("B20").Value = 0.7018145335602526
$ws.Range("C20").Value = 0.296059239290754
$ws.Range("D20").Value = 0.04634991759733964
$ws.Range("F20").Value = 0.7003488478517639
$ws.Range("G20").Value = 0.5355114058765125
$ws.Range("H20").Value = 0.6723908490769617
$ws.Range("K20").Value = 0.4143746804722923
$ws.Range("L20").Value = 0.3437314048043731
$ws.Range("O20").Value = 2.385981861512533
$ws.Range("B21").Value = 0.7689959453883262
$ws.Range("C21").Value = 0.2926889478837964
$ws.Range("D21").Value = 0.05013700497318041
$ws.Range("F21").Value = 0.701413768962162
$ws.Range("G21").Value = 0.5338794512415177
$ws.Range("H21").Value = 0.6674430866165864
$ws.Range("K21").Value = 0.4651495014642819
$ws.Range("L21").Value = 0.3610105313767207
$ws.Range("O21").Value = 2.372443647231023
$ws.Range("B22").Value = 0.8129522258245174
$ws.Range("C22").Value = 0.2905750458513747
$ws.Range("D22").Value = 0.0525987637908969
$ws.Range("F22").Value = 0.7024933558345765
$ws.Range("G22").Value = 0.5331441462222983
$ws.Range("H22").Value = 0.6644694471450663
$ws.Range("K22").Value = 0.49824733415889
$ws.Range("L22").Value = 0.3724542471237271
$ws.Range("O22").Value = 2.364834786595367
$ws.Range("B23").Value = 0.7894876905600654
$ws.Range("C23").Value = 0.2916951580322209
$ws.Range("D23").Value = 0.05128610967572911
$ws.Range("F23").Value = 0.7018820612428769
$ws.Range("G23").Value = 0.5335062533460473
$ws.Range("H23").Value = 0.6660327539721749
$ws.Range("K23").Value = 0.4805905179695458
$ws.Range("L23").Value = 0.366332749553834
$ws.Range("O23").Value = 2.368782358213508
$ws.Range("B24").Value = 0.700717382366804
$ws.Range("C24").Value = 0.2961158250950433
$ws.Range("D24").Value = 0.04628779783225667
$ws.Range("F24").Value = 0.7003379183171816
$ws.Range("G24").Value = 0.5355436656845143
$ws.Range("H24").Value = 0.6724761095024405
$ws.Range("K24").Value = 0.4135433798317365
$ws.Range("L24").Value = 0.3434515385484502
$ws.Range("O24").Value = 2.386224047262857
$ws.Range("B25").Value = 0.6052034198995386
$ws.Range("C25").Value = 0.3012660273729626
$ws.Range("D25").Value = 0.04084050487688273
$ws.Range("F25").Value = 0.7003192319041887
$ws.Range("G25").Value = 0.5391607787195056
$ws.Range("H25").Value = 0.6805436819168591
$ws.Range("K25").Value = 0.3408720412251114
$ws.Range("L25").Value = 0.3194205941220645
$ws.Range("O25").Value = 2.410353852112578
